$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate formatting of the last BOM row (21) onto the new row (22)
$ws.Range("A21:G21").Copy()
$ws.Range("A22:G22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the new "Axle boots" BOM line (order chosen to match
# the shared-string insertion order: ID, Comments, Part)
$ws.Range("G22").Value = "EN_07020"
$ws.Range("E22").Value = "Over driveshafts and tripod housings"
$ws.Range("C22").Value = "Axle boots"
$ws.Range("D22").Value = "b"
$ws.Range("F22").Value = 4

$ws.Range("C23").Select()
